# Apply the dated-output refresh: update the header date and all
# three-digit-by-one-digit multiplication problems/answers in the table.

$d = $word.ActiveDocument

# Mapping of old text -> new text (each old value occurs exactly once).
$replacements = [ordered]@{
    "2026-02-19 Thursday" = "2026-02-20 Friday"
    "429×7=3003" = "636×2=1272"
    "114×3=342"  = "706×8=5648"
    "723×9=6507" = "665×4=2660"
    "922×8=7376" = "744×3=2232"
    "278×3=834"  = "612×3=1836"
    "336×6=2016" = "686×6=4116"
    "132×9=1188" = "728×3=2184"
    "550×8=4400" = "626×3=1878"
    "936×5=4680" = "907×4=3628"
    "146×8=1168" = "712×5=3560"
    "684×3=2052" = "808×4=3232"
    "278×8=2224" = "588×6=3528"
    "625×2=1250" = "252×8=2016"
    "361×7=2527" = "233×6=1398"
    "743×7=5201" = "568×6=3408"
    "450×3=1350" = "196×2=392"
    "807×2=1614" = "590×5=2950"
    "888×4=3552" = "540×5=2700"
    "590×3=1770" = "283×8=2264"
    "488×2=976"  = "591×7=4137"
    "372×5=1860" = "913×8=7304"
    "801×5=4005" = "700×8=5600"
    "378×8=3024" = "609×4=2436"
    "581×3=1743" = "461×3=1383"
    "607×4=2428" = "737×4=2948"
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
